$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "291.22"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-3.19%"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-4.94%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.954"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.33%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07217"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-5.77%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.817"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-5.75%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.699"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-1.82%"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.78%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8973"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-2.09%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1658"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-5.58%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07704"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-0.55%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08016"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-6.08%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03041"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-4.07%"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.11%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001506"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.20%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005698"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-3.74%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.468"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.07%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.083"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-3.29%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3317"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-1.02%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1331"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.40%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.052"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-5.45%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.2389"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "19.99%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.04498"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.35%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001215"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.62%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-8.85%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001251"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "0.02%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01581"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-7.17%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04400"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-6.43%"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-2.67%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01006"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1305"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-3.27%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002062"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-11.57%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009191"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-12.38%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00005950"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-4.76%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.02%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.247"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "173.92%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002102"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.02%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002001"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.02%"
